$wb = $excel.ActiveWorkbook

# --- 1er Parcial ---
$ws1 = $wb.Worksheets.Item("1er Parcial")
$ws1.Range("I4").Value = 6.2

# --- 2o Parcial ---
$ws2 = $wb.Worksheets.Item("2o Parcial")
$ws2.Range("E2").Value = 36
$ws2.Range("F2").Value = 1
$ws2.Range("G2").Value = 97.3
$ws2.Range("H2").Value = 2.7
$ws2.Range("I2").Value = 7.3
$ws2.Range("J2").Value = 0
$ws2.Range("K2").Value = 0
$ws2.Range("E3").Value = 39
$ws2.Range("F3").Value = 0
$ws2.Range("G3").Value = 100
$ws2.Range("H3").Value = 0
$ws2.Range("I3").Value = 8.2
$ws2.Range("J3").Value = 0
$ws2.Range("K3").Value = 0
$ws2.Range("E4").Value = 29
$ws2.Range("F4").Value = 3
$ws2.Range("G4").Value = 90.6
$ws2.Range("H4").Value = 9.4
$ws2.Range("I4").Value = 7.3
$ws2.Range("J4").Value = 0
$ws2.Range("K4").Value = 0
$ws2.Range("E5").Value = 47
$ws2.Range("F5").Value = 2
$ws2.Range("G5").Value = 95.9
$ws2.Range("H5").Value = 4.1
$ws2.Range("I5").Value = 7.8
$ws2.Range("J5").Value = 0
$ws2.Range("K5").Value = 0
$ws2.Range("E6").Value = 38
$ws2.Range("F6").Value = 1
$ws2.Range("G6").Value = 97.4
$ws2.Range("H6").Value = 2.6
$ws2.Range("I6").Value = 6.8
$ws2.Range("J6").Value = 0
$ws2.Range("K6").Value = 0
$ws2.Range("E7").Value = 36
$ws2.Range("F7").Value = 2
$ws2.Range("G7").Value = 94.7
$ws2.Range("H7").Value = 5.3
$ws2.Range("I7").Value = 7.9
$ws2.Range("J7").Value = 0
$ws2.Range("K7").Value = 0
$ws2.Range("E8").Value = 225
$ws2.Range("F8").Value = 9
$ws2.Range("G8").Value = 96.2
$ws2.Range("H8").Value = 3.8
$ws2.Range("I8").Value = 7.6
$ws2.Range("J8").Value = 0
$ws2.Range("K8").Value = 0
$ws2.Range("E15").Value = 225
$ws2.Range("F15").Value = 137
$ws2.Range("G15").Value = 62.2
$ws2.Range("H15").Value = 37.8
$ws2.Range("I15").Value = 4.1
$ws2.Range("J15").Value = 128
$ws2.Range("K15").Value = 35.4

# --- Final ---
$ws3 = $wb.Worksheets.Item("Final")
$ws3.Range("I2").Value = 7
$ws3.Range("E3").Value = 39
$ws3.Range("F3").Value = 0
$ws3.Range("G3").Value = 100
$ws3.Range("H3").Value = 0
$ws3.Range("I3").Value = 7.4
$ws3.Range("E4").Value = 29
$ws3.Range("F4").Value = 3
$ws3.Range("G4").Value = 90.6
$ws3.Range("H4").Value = 9.4
$ws3.Range("I4").Value = 6.8
$ws3.Range("I5").Value = 7.2
$ws3.Range("I6").Value = 6.5
$ws3.Range("E7").Value = 36
$ws3.Range("F7").Value = 2
$ws3.Range("G7").Value = 94.7
$ws3.Range("H7").Value = 5.3
$ws3.Range("I7").Value = 7.3
$ws3.Range("E8").Value = 225
$ws3.Range("F8").Value = 9
$ws3.Range("G8").Value = 96.2
$ws3.Range("H8").Value = 3.8
$ws3.Range("I8").Value = 7
$ws3.Range("E15").Value = 286
$ws3.Range("F15").Value = 76
$ws3.Range("G15").Value = 79
$ws3.Range("H15").Value = 21
$ws3.Range("I15").Value = 6.7
